$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19 (Vimex slide switch SS-12E17) was only partially filled in before -
# add the source link plus the price/qty/shipping figures now that the part
# has actually been bought & the gain-scheduled controller is working.
$ws.Range("E19").Formula = "=2.94/6"
$ws.Range("F19").Value = 1
$ws.Range("G19").Formula = "=E19*F19"
$ws.Range("H19").Formula = "=(3.52+5)/4"
$ws.Range("I19").Formula = "=G19+H19"

$url = "http://www.jameco.com/z/SS-12E17-3-Pin-SPDT-Slide-Switch_2258831.html"
$ws.Hyperlinks.Add($ws.Range("D19"), $url)
$ws.Range("D19").Borders.LineStyle = 1

# Leave the view where the user was last working
$ws.Activate()
$ws.Range("B19").Select()
